$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3166
$ws1.Range("F5").Value = 6862
$ws1.Range("F6").Value = 1964
$ws1.Range("F12").Value = 18
$ws1.Range("F14").Value = 174

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3166
$ws4.Range("F6").Value = 6862
$ws4.Range("F7").Value = 1964
$ws4.Range("F13").Value = 18
$ws4.Range("F15").Value = 174
